$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format before writing, so numeric-looking
# strings (e.g. "1.004", "0.9992") are stored as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.746.71'
$ws.Range("E2").Value = '  +5.37%  '
$ws.Range("D3").Value = '1.702.34'
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '330.41'
$ws.Range("E5").Value = '  +6.02%  '
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.3683'
$ws.Range("E7").Value = '  +0.93%  '
$ws.Range("D8").Value = '48.48'
$ws.Range("E8").Value = '  +4.26%  '
$ws.Range("D9").Value = '0.3303'
$ws.Range("E9").Value = '  +1.96%  '
$ws.Range("E10").Value = '  +4.01%  '
$ws.Range("D11").Value = '0.07319'
$ws.Range("E11").Value = '  +4.42%  '
$ws.Range("D12").Value = '0.9992'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '6.182'
$ws.Range("E13").Value = '  +3.81%  '
$ws.Range("D14").Value = '19.97'
$ws.Range("E14").Value = '  +3.25%  '
$ws.Range("D15").Value = '6.845'
$ws.Range("E15").Value = '  +3.81%  '
$ws.Range("D16").Value = '1.702.79'
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").Value = '0.00001063'
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").Value = '0.06625'
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("D19").Value = '81.05'
$ws.Range("E19").Value = '  +3.42%  '
$ws.Range("D21").Value = '16.12'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").Value = '6.029'
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").Value = '12.94'
$ws.Range("E23").Value = '  +3.36%  '
$ws.Range("D24").Value = '25.743.33'
$ws.Range("E24").Value = '  +5.38%  '
$ws.Range("D25").Value = '2.458'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '2.472'
$ws.Range("E26").Value = '  +6.70%  '
$ws.Range("D27").Value = '149.47'
$ws.Range("E27").Value = '  +2.08%  '
$ws.Range("D28").Value = '19.16'
$ws.Range("E28").Value = '  +3.36%  '
$ws.Range("D29").Value = '1.286'
$ws.Range("E29").Value = '  +8.51%  '
$ws.Range("D30").Value = '1.891.30'
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").Value = '128.00'
$ws.Range("E31").Value = '  +3.43%  '
$ws.Range("D32").Value = '4.097'
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").Value = '5.913'
$ws.Range("E33").Value = '  +4.08%  '
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.712'
$ws.Range("E34").Value = '  +3.75%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08474'
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").Value = '12.81'
$ws.Range("E36").Value = '  +6.17%  '
$ws.Range("D37").Value = '5.307'
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("D38").Value = '1.276'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '0.06168'
$ws.Range("E39").Value = '  +2.64%  '
$ws.Range("D40").Value = '8.492'
$ws.Range("E40").Value = '  +5.23%  '
$ws.Range("D41").Value = '0.2111'
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("D42").Value = '0.02240'
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("D43").Value = '14.68'
$ws.Range("E43").Value = '  +17.28%  '
$ws.Range("D44").Value = '0.6091'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("E47").Value = '  +3.77%  '
$ws.Range("D48").Value = '126.29'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("D49").Value = '1.994'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").Value = '0.07211'
$ws.Range("E50").Value = '  +4.50%  '
$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
$ws.Range("D51").Value = '1.203'
$ws.Range("E51").Value = '  +2.15%  '

# Restore the default (General) style so no residual formatting is left behind.
$ws.Range("D2:E51").Style = "Normal"
